# Regenerate the localization-status report: the handback pass has moved
# every file that was "Ready for handoff" into "In Translation", so that
# status string is no longer used anywhere in the workbook.
#
# Overview sheet: zh-cn / de-de status columns (E, F)
# zh-cn / de-de sheets: Status column (C)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E5").Value = "In Translation"
$ws.Range("F5").Value = "In Translation"
$ws.Range("E6").Value = "In Translation"
$ws.Range("F6").Value = "In Translation"
$ws.Range("E7").Value = "In Translation"
$ws.Range("F7").Value = "In Translation"
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C5").Value = "In Translation"
$wsZh.Range("C6").Value = "In Translation"
$wsZh.Range("C7").Value = "In Translation"
$wsZh.Columns.Item(3).AutoFit()

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C5").Value = "In Translation"
$wsDe.Range("C6").Value = "In Translation"
$wsDe.Range("C7").Value = "In Translation"
$wsDe.Columns.Item(3).AutoFit()
